$d = $word.ActiveDocument

# --- 1. Remove the _GoBack bookmark from its current (empty) paragraph ---
# It will be re-added later, attached to its new location.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Locate the paragraph that ends with "...caso seja finalizado." ---
# (the last paragraph of the document) and append the new content after it.
$last = $d.Paragraphs.Last
$r = $last.Range
$r.Collapse(0)

# blank paragraph
$r.InsertAfter("`r")
$r.Collapse(0)

# "A PlayList de Musicas..." paragraph
$r.InsertAfter("`r")
$r.Collapse(0)
$r.InsertAfter('A PlayList de Músicas (#musica, #artista, album, nomeArquivo, endereco), Inserção de Vinhetas e Spots (#musica, #artista, album, nomeArquivo, endereco), Programações Gravadas (#musica, #artista, album, nomeArquivo, endereco), todas são referencias a arquivos de áudio padrão e independentes entre sí.')
$r.Collapse(0)

# blank paragraph
$r.InsertAfter("`r")
$r.Collapse(0)

# "O Gerenciamento de Eventos Automatizados..." paragraph
$r.InsertAfter("`r")
$r.Collapse(0)
$r.InsertAfter('O Gerenciamento de Eventos Automatizados (#data, #hora, #minuto, #diaSemana, nomeEvento, prioridade) possui como referencia um único Arquivo de Áudio (#musica, #artista, album, nomeArquivo, endereco).')
$r.Collapse(0)

# blank paragraph
$r.InsertAfter("`r")
$r.Collapse(0)

# "O Gerenciamento de Eventos Externos Automatizados..." paragraph (bookmark sits inside it)
$r.InsertAfter("`r")
$r.Collapse(0)
$r.InsertAfter('O Gerenciamento de Eventos Externos Automatizados (#data, #hora, #minuto, #diaSemana, nomeEvento, prioridade) possui como referencia um único Streaming de Áudio (#Evento')
$r.Collapse(0)
$r.InsertAfter(', endereco).')
$r.Collapse(0)

# two blank paragraphs
$r.InsertAfter("`r")
$r.Collapse(0)
$r.InsertAfter("`r")
$r.Collapse(0)

# "O Gerenciamento de Eventos Complexos..." paragraph
$r.InsertAfter("`r")
$r.Collapse(0)
$r.InsertAfter('O Gerenciamento de Eventos Complexos (#data, #hora, #minuto, #diaSemana, nomeEvento, prioridade) possui como referencia Comandos Específicos (#comando, #parametro, parametroOpcional).')
$r.Collapse(0)

# --- 3. Re-attach the _GoBack bookmark inside the "Streaming de Audio" paragraph, ---
# right after "...(#Evento" and before ", endereco)."
$marker = '(#Evento'
$searchRange = $d.Content
$found = $searchRange.Find.Execute($marker, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $bmPos = $searchRange.End
    $bmRange = $d.Range($bmPos, $bmPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

Write-Output "Paragraphs: $($d.Paragraphs.Count)"
Write-Output "Bookmark exists: $($d.Bookmarks.Exists('_GoBack'))"
